$d = $word.ActiveDocument

# Use Track Changes so each InsertAfter() call is recorded as its own
# insertion (and therefore serialized as its own <w:r> run) instead of
# being silently coalesced into the run it's adjacent to.
$d.TrackRevisions = $true

$para = $d.Paragraphs.Item(1)
$r = $para.Range

$r.InsertAfter(" (")
$r.InsertAfter("Changed main")
$r.InsertAfter(")")

$d.TrackRevisions = $false

# Accept the tracked insertions one at a time (rather than a single
# AcceptAllRevisions() sweep) so the rest of the document's cached
# layout info is left untouched.
while ($d.Revisions.Count -gt 0) {
  $d.Revisions.Item(1).Accept()
}
